# Update quizvragen via Admin
# Adds two new multiple-choice question rows to the "DC" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DC")

# Row 3 - new question: "Wat betekend gelijkstroom"
$ws.Range("B3").Value = "mc"
$ws.Range("D3").Value = "Wat betekend gelijkstroom"
$ws.Range("E3").Value = "['test', 'test1', 'test2']"
$ws.Range("F3").Value = 1

# Row 4 - new question: "Wat betekend gelijkstroom1"
$ws.Range("B4").Value = "mc"
$ws.Range("D4").Value = "Wat betekend gelijkstroom1"
$ws.Range("E4").Value = "['test', 'test1', 'test2']"
$ws.Range("F4").Value = 1
